# Shift each sheet's date/amount series up by two rows (dropping the two
# oldest dates) and append two new trailing rows (dates 45959, 45960) with
# a placeholder value of 0, matching the refreshed daily export.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Pull the old rows 4:101 (A:B) as a 2D array, then write them into rows 2:99.
    $src = $ws.Range("A4:B101")
    $vals = $src.Value2
    $dst = $ws.Range("A2:B99")
    $dst.Value2 = $vals

    # Append the two new trailing rows with placeholder zero amounts.
    $ws.Range("A100").Value2 = 45959
    $ws.Range("B100").Value2 = 0
    $ws.Range("A101").Value2 = 45960
    $ws.Range("B101").Value2 = 0
}
